$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 2.964545797025059)
  3  = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 3.536033448013082)
  4  = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 17.65757632934944)
  5  = @(1.445647641019636, 1.626987699542094, 18.71679738969934, 13.86384647080068, 0, 35.65327920106175)
  6  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 6.82939032824165)
  7  = @(3.272327238179451, 1.626987699542094, 189.6080260415259, 13.86384647080068, 1, 208.3711874500482)
  8  = @(1.445647641019636, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 22.32281868886277)
  9  = @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 1.84748871573303)
  10 = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 6.038307959104277)
  11 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 0, 24.14949828602258)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  $ws.Cells.Item($r, 2).Value = $vals[0]
  $ws.Cells.Item($r, 3).Value = $vals[1]
  $ws.Cells.Item($r, 4).Value = $vals[2]
  $ws.Cells.Item($r, 5).Value = $vals[3]
  $ws.Cells.Item($r, 6).Value = $vals[4]
  $ws.Cells.Item($r, 7).Value = $vals[5]
}
